$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account")

$ws.Range("C1").Value = "name"
$ws.Range("E1").Value = "password"
$ws.Range("F1").Value = "consent"

$ws.Range("F2").Value = "I have consented"
$ws.Range("F3").Value = "I have consented"
$ws.Range("F4").Value = "I have consented"
$ws.Range("F5").Value = "I have consented"
$ws.Range("F6").Value = "I have consented"
$ws.Range("F7").Value = "I have consented"
$ws.Range("F8").Value = "I have consented"
$ws.Range("F9").Value = "I have consented"
$ws.Range("F10").Value = "I have consented"

$ws.Range("C2").Value = "Tom Yeh"
$ws.Range("C3").Value = "Abby Stangl"
$ws.Range("C4").Value = "Mike Skirpan"
$ws.Range("C5").Value = "Jose Meti"
$ws.Range("C6").Value = "Matt Kesh"
$ws.Range("C7").Value = "Jackie Hama"
$ws.Range("C8").Value = "Carol Boston"
$ws.Range("C9").Value = "Jenny Preece"
$ws.Range("C10").Value = "Petter Joisterest"

$ws.Range("E2").Value = 1234
$ws.Range("E3").Value = 1234
$ws.Range("E4").Value = 1234
$ws.Range("E5").Value = 1234
$ws.Range("E6").Value = 1234
$ws.Range("E7").Value = 1234
$ws.Range("E8").Value = 1234
$ws.Range("E9").Value = 1234
$ws.Range("E10").Value = 1234

$ws.Range("D1").Value = "email"

$ws.Range("D2").Value = "tom@mail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:tom@mail.com")
$ws.Range("D3").Value = "abby@mail.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:abby@mail.com")
$ws.Range("D4").Value = "mike@mail.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:mike@mail.com")
$ws.Range("D5").Value = "jose@mail.com"
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:jose@mail.com")
$ws.Range("D6").Value = "matt@mail.com"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:matt@mail.com")
$ws.Range("D7").Value = "jackie@mail.com"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:jackie@mail.com")
$ws.Range("D8").Value = "carol@umd.edu"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:carol@umd.edu")
$ws.Range("D9").Value = "jenny@umd.edu"
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:jenny@umd.edu")
$ws.Range("D10").Value = "petter@umd.edu"
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:petter@umd.edu")

$ws.Columns.Item(3).ColumnWidth = 8.1666666666667
$ws.Columns.Item(4).ColumnWidth = 14.1666666666667
$ws.Columns.Item(5).ColumnWidth = 13.6666666666667

$ws.Activate()
$ws.Range("E10").Select()

Write-Output "done"
